$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of trade data appended below the existing history (row 9)
$ws.Range("A9").Value = 9895.27
$ws.Range("B9").Value = 9974.07
$ws.Range("C9").Value = 282.89999999999998
$ws.Range("D9").Value = 285.14
$ws.Range("E9").Value = $true
$ws.Range("F9").Value = 0.79

# Copy the date/time number formatting from the row above so G9 reuses
# the existing style (rather than minting a new one), then set its value.
$ws.Range("G8").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = 42609.487268518518

$ws.Range("H9").Value = $false
